# fix: Improved error message during download
# Relocates the "error downloading video" message (row 31) to include
# extra guidance for the user, and updates the saved view state
# (scroll position / selection) to match where the author left off.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Update the download-error message (row 31) ------------------------
# Spanish (column A) and English (column B) versions both gain two extra
# lines inviting the user to get in touch if the problem persists.
$ws.Range("A31").Value = " Ocurrió un error al descargar el video. `nInténtalo nuevamente más tarde.`nSi el error persiste, ponte en contacto conmigo."
$ws.Range("B31").Value = " An error occurred while downloading the video. `nPlease try again later.`nIf the error persists, please contact me."

# Match the formatting used by the other multi-line messages in this
# sheet (rows 24-26): centered, wrapped text, auto height for 3 lines.
$errRow = $ws.Range("A31:B31")
$errRow.HorizontalAlignment = -4108
$errRow.WrapText = $true
$ws.Rows.Item(31).RowHeight = 45

# --- Restore the workbook view / selection seen in the edited file -----
$ws.Activate()
$excel.ActiveWindow.ScrollRow = 17
$excel.ActiveWindow.ScrollColumn = 1
$ws.Range("D27").Select()
